# Sprint 7 update: insert a new "Sprint 7" status slide between the
# existing "Sprint 7" slide and the "Próxima Sprint" slide, describing
# the tests added and the Arduino/Bluetooth/Android communication work.

$p = $ppt.ActivePresentation

# Insert a new Title+Content slide at position 2 (pushes "Próxima Sprint"
# down to position 3). Layout 2 = "Título e conteúdo", same layout used
# by the other slides in this deck.
$s = $p.Slides.Add(2, 2)

# --- Title placeholder: "Sprint 7" (centered, like the other slides) ---
$title = $s.Shapes.Item(1).TextFrame.TextRange
$title.Text = "Sprint 7"
$title.LanguageID = "pt-BR"
$title.ParagraphFormat.Alignment = 2

# --- Content placeholder ---
$body = $s.Shapes.Item(2).TextFrame.TextRange

# Paragraph 1
$body.Text = "Adicionamos os testes de para quando os valores forem maior que o desejado."
$body.LanguageID = "pt-BR"

# Paragraph 2 - built run by run so wording stays in separate runs
$r = $body.InsertAfter("`rRealizamos a comunicação com ")
$r.LanguageID = "pt-BR"
$r = $r.InsertAfter("arduino")
$r.LanguageID = "pt-BR"
$r = $r.InsertAfter(" ")
$r.LanguageID = "pt-BR"
$r = $r.InsertAfter("via ")
$r.LanguageID = "pt-BR"
$r = $r.InsertAfter("bluetooth")
$r.LanguageID = "pt-BR"
$r = $r.InsertAfter(" com o ")
$r.LanguageID = "pt-BR"
$r = $r.InsertAfter("android")
$r.LanguageID = "pt-BR"
$r = $r.InsertAfter(".")
$r.LanguageID = "pt-BR"

# Paragraph 3 - trailing empty paragraph, no bullet
$r = $r.InsertAfter("`r")

$lastPara = $body.Paragraphs(3, 1)
$lastPara.ParagraphFormat.Bullet.Visible = $false
